$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# New benchmark "book": row 14 = base run, row 15 = our run (mirrors the
# existing web/base/our block layout used by every other named benchmark).
# ---------------------------------------------------------------------------

# --- Row 14 (book / base) ---
$ws.Range("A14").Value = "book"
$ws.Range("B14").Value = "base"

$ws.Range("C3").Copy()
$ws.Range("C14").PasteSpecial($xlPasteFormats)
$ws.Range("C14").Value = 120.212642

$ws.Range("D3").Copy()
$ws.Range("D14").PasteSpecial($xlPasteFormats)
$ws.Range("D14").Value = 0.160019

$ws.Range("E14").Formula = "=C14*D14"

$ws.Range("F3").Copy()
$ws.Range("F14").PasteSpecial($xlPasteFormats)
$ws.Range("F14").Value = 3.723402

$ws.Range("G3").Copy()
$ws.Range("G14").PasteSpecial($xlPasteFormats)
$ws.Range("G14").Value = 0.893672

# --- Row 15 (book / our) ---
$ws.Range("B15").Value = "our"

$ws.Range("C3").Copy()
$ws.Range("C15").PasteSpecial($xlPasteFormats)
$ws.Range("C15").Value = 120.96821

$ws.Range("D3").Copy()
$ws.Range("D15").PasteSpecial($xlPasteFormats)
$ws.Range("D15").Value = 0.164621

$ws.Range("E15").Formula = "=C15*D15"

$ws.Range("F3").Copy()
$ws.Range("F15").PasteSpecial($xlPasteFormats)
$ws.Range("F15").Value = 1.771538

$ws.Range("G3").Copy()
$ws.Range("G15").PasteSpecial($xlPasteFormats)
$ws.Range("G15").Value = 47.31291

$excel.CutCopyMode = $false

# --- Window/view state: zoom in and leave selection where editing stopped ---
$ws.Range("D16").Select() | Out-Null
$excel.ActiveWindow.Zoom = 175
